$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37 (pushes old rows 37..56 down to 38..57),
# then fill in the new weekly record (week of 2022-01-06).
$ws.Rows("37:37").Insert()

$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44567
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 100112031
$ws.Cells.Item(37, 7).Value = "Poroto verde"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 120
$ws.Cells.Item(37, 11).Value = 38000
$ws.Cells.Item(37, 12).Value = 39000
$ws.Cells.Item(37, 13).Value = 38500
$ws.Cells.Item(37, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Región del Maule"
$ws.Cells.Item(37, 16).Value = 1540
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"

# Insert a second new row at 56 (after the first insert, the record
# that used to be at old row 55 - 2021-12-10, Magnum - now sits at row
# 56; push it down to row 57, where it belongs unchanged, and fill the
# freed row 56 with the other new weekly record, week of 2022-01-07).
$ws.Rows("56:56").Insert()

$ws.Cells.Item(56, 1).Value = 7
$ws.Cells.Item(56, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(56, 3).Value = "Ñuble"
$ws.Cells.Item(56, 4).Value = 44568
$ws.Cells.Item(56, 5).Value = 16
$ws.Cells.Item(56, 6).Value = 100112031
$ws.Cells.Item(56, 7).Value = "Poroto verde"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 120
$ws.Cells.Item(56, 11).Value = 37000
$ws.Cells.Item(56, 12).Value = 38000
$ws.Cells.Item(56, 13).Value = 37500
$ws.Cells.Item(56, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Región del Maule"
$ws.Cells.Item(56, 16).Value = 1500
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"

# Rows 57 (2021-12-10, Magnum) and 58 (2021-03-17, Región de O'Higgins)
# now already hold the correct shifted-down values and need no
# further edits.
